# Penyesuaian dengan database: pergantian kolom id_x dengan nama_x di tabel dokumen.
# Di sheet "Lembar1", kolom B (sebelumnya berisi kode id "001") diganti
# menjadi nama/kode "010" untuk baris data 2-4, dan posisi seleksi aktif
# dipindahkan ke G5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lembar1")

$ws.Range("B2").Value = "010"
$ws.Range("B3").Value = "010"
$ws.Range("B4").Value = "010"

$ws.Activate()
$ws.Range("G5").Select()
